$wb = $excel.ActiveWorkbook

# --- Simulated Wild Card round: append new game play-by-play data
# --- to the running season lists (YDS rush/pass yards, ST return
# --- yards/distances), then update the derived season totals. ---

$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = $ws.Range("B2").Value() + " 1 1 5 9 -1 2 2 -6 1 1 2 0 7 14 1 -1 3 17 -1 0 2 4 3 0 2 0 15 2 -2 2 5 2 4 3 4 4 8 0 -2 2 4 7 2 13 4"
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("C2").Value = $ws.Range("C2").Value() + " 2 2 0 4 22 -5 2 3 3 4 2 3 18 8 -1 8 10 46 1 22 27 3 6 0 4 2 6 1 21 0 11 4 4 3 1 6 4 4 4 2 2 -1 23 -2 15 13 12 1 0 -2 -4 5 4 3 6 11 -1"
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B3").Value = $ws.Range("B3").Value() + " 4 11 10 7 0 7 25 2 5 16 8 2 10 10 7 8 -2 8 20 11 3 6 10 6 3 11 14 1 11 3 4 7 5 5 0 16 0 6 8 5 5 22 13 1 4 7 9 8 15 5 5 7 7 10 12 11 11"
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("C3").Value = $ws.Range("C3").Value() + " 3 19 13 11 10 6 0 28 7 1 5 1 6 11 8 12 -4 20 5 5 4 15 20 2 14 15 5 4 31 0 3 12 15 27 48 29 8 23 4 1 31 6 15 3 41 2 2"
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B4").Value = $ws.Range("B4").Value() + " 53 56 58 59 59 64 60"
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B5").Value = $ws.Range("B5").Value() + " 19 23 27 26 18 23 27"
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B6").Value = $ws.Range("B6").Value() + " 32 13 21 23 14 29"
$ws = $wb.Worksheets.Item("ST")
$ws.Range("D3").Value = $ws.Range("D3").Value() + " 36 34 26 49 33 46 36 38 54 56 55 42 47 41 53"
$ws = $wb.Worksheets.Item("ST")
$ws.Range("D4").Value = $ws.Range("D4").Value() + " 0 0 0 0 0 0 0 0 21 1 48 0 0 0 -8"
$ws = $wb.Worksheets.Item("ST")
$ws.Range("D5").Value = $ws.Range("D5").Value() + " 7 0 0 0 23 0 0 8 0 0 5"

# --- Update season totals rows on each summary sheet ---

$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 372
$ws.Range("F2").Value = 135
$ws.Range("G2").Value = 93
$ws.Range("J2").Value = 57
$ws.Range("N2").Value = 30
$ws.Range("O2").Value = 44
$ws.Range("P2").Value = 25
$ws.Range("B3").Value = 17
$ws.Range("C3").Value = 463
$ws.Range("E3").Value = 68
$ws.Range("F3").Value = 271
$ws.Range("G3").Value = 114
$ws.Range("H3").Value = 65
$ws.Range("I3").Value = 168
$ws.Range("J3").Value = 166
$ws.Range("L3").Value = 686
$ws.Range("M3").Value = 437
$ws.Range("Q3").Value = 1098
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 473
$ws.Range("D2").Value = 30
$ws.Range("F2").Value = 128
$ws.Range("G2").Value = 108
$ws.Range("I2").Value = 20
$ws.Range("J2").Value = 65
$ws.Range("N2").Value = 42
$ws.Range("O2").Value = 44
$ws.Range("P2").Value = 18
$ws.Range("B3").Value = 21
$ws.Range("C3").Value = 338
$ws.Range("D3").Value = 11
$ws.Range("E3").Value = 90
$ws.Range("F3").Value = 235
$ws.Range("G3").Value = 67
$ws.Range("H3").Value = 80
$ws.Range("I3").Value = 134
$ws.Range("J3").Value = 106
$ws.Range("L3").Value = 573
$ws.Range("M3").Value = 346
$ws.Range("Q3").Value = 1125
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 182
$ws.Range("D2").Value = 163
$ws.Range("F2").Value = 248
$ws.Range("G2").Value = 234
$ws.Range("J2").Value = 110
$ws.Range("K2").Value = 104
$ws.Range("L2").Value = 75
$ws.Range("M2").Value = 61
$ws.Range("B3").Value = 88
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("B3").Value = 16
$ws.Range("C3").Value = 17
$ws.Range("D3").Value = 17
$ws.Range("E3").Value = 17
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 35
$ws.Range("B3").Value = 33
